# Insert a new daily price record as row 114, pushing the existing rows
# (114-194) down to (115-195), matching the "Fruta / hortaliza, semanal"
# update for Arveja Verde / Mercado Mayorista Lo Valledor de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 114:194 down to 115:195 by inserting a new blank row at 114.
$ws.Rows.Item(114).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(114, 1).Value = 6
$ws.Cells.Item(114, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(114, 3).Value = "Metropolitana"
$ws.Cells.Item(114, 4).Value = 44589
$ws.Cells.Item(114, 5).Value = 13
$ws.Cells.Item(114, 6).Value = 100112022
$ws.Cells.Item(114, 7).Value = "Arveja Verde"
$ws.Cells.Item(114, 8).Value = "Perfection"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 230
$ws.Cells.Item(114, 11).Value = 20000
$ws.Cells.Item(114, 12).Value = 23000
$ws.Cells.Item(114, 13).Value = 21696
$ws.Cells.Item(114, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(114, 15).Value = "Carahue"
$ws.Cells.Item(114, 16).Value = 868
$ws.Cells.Item(114, 17).Value = 25
$ws.Cells.Item(114, 18).Value = "Hortaliza"
